$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B3").Value = 5.785
$ws.Range("A12").Value = -21.626
$ws.Range("B14").Value = 5.944
$ws.Range("B26").Value = 6.706
$ws.Range("B31").Value = 6.205
$ws.Range("A32").Value = -21.351
$ws.Range("B35").Value = 8.095000000000001
$ws.Range("A36").Value = -20.945
$ws.Range("B37").Value = 8.260000000000002
$ws.Range("A38").Value = -20.093
$ws.Range("B45").Value = 5.786
$ws.Range("A46").Value = -21.553
$ws.Range("A54").Value = -22.209
$ws.Range("A55").Value = -22.164
$ws.Range("B57").Value = 5.331999999999999
$ws.Range("A67").Value = -21.588
$ws.Range("A69").Value = -21.721
$ws.Range("A72").Value = -21.55
$ws.Range("A91").Value = -21.522
$ws.Range("A99").Value = -20.828
$ws.Range("B100").Value = 5.558
$ws.Range("B102").Value = 7.186
